$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds text values (e.g. "3.137.72") that Excel would
# otherwise auto-convert to numbers; force it to stay text before writing.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "64.705.67"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "3.139.48"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "572.18"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").Value = "147.95"
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("D8").Value = "3.137.44"
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  -0.91%  "
$ws.Range("E10").Value = "  -4.20%  "
$ws.Range("D11").Value = "6.06"
$ws.Range("E11").Value = "  -1.71%  "
$ws.Range("E12").Value = "  -1.13%  "
$ws.Range("E13").Value = "  -1.66%  "
$ws.Range("D14").Value = "36.82"
$ws.Range("E14").Value = "  -0.81%  "
$ws.Range("D15").Value = "3.651.67"
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("D16").Value = "64.834.40"
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("D17").Value = "3.139.03"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").Value = "7.05"
$ws.Range("E18").Value = "  -1.35%  "
$ws.Range("E19").Value = "  -0.22%  "
$ws.Range("D20").Value = "499.40"
$ws.Range("E20").Value = "  -2.43%  "
$ws.Range("D21").Value = "14.72"
$ws.Range("E21").Value = "  -0.72%  "
$ws.Range("E22").Value = "  -2.28%  "
$ws.Range("D23").Value = "15.07"
$ws.Range("E23").Value = "  -2.81%  "
$ws.Range("E24").Value = "  -1.99%  "
$ws.Range("D25").Value = "83.53"
$ws.Range("E25").Value = "  -2.25%  "
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("E27").Value = "  -1.50%  "
$ws.Range("D28").Value = "8.77"
$ws.Range("E28").Value = "  +0.35%  "
$ws.Range("E29").Value = "  +0.92%  "
$ws.Range("B30").Value = "Stacks"
$ws.Range("C30").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D30").Value = "2.76"
$ws.Range("E30").Value = "  +3.54%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "27.31"
$ws.Range("E31").Value = "  -2.37%  "
$ws.Range("E32").Value = "  -0.22%  "
$ws.Range("E33").Value = "  -0.61%  "
$ws.Range("D34").Value = "6.11"
$ws.Range("E34").Value = "  +1.20%  "
$ws.Range("D35").Value = "6.41"
$ws.Range("E35").Value = "  -2.63%  "
$ws.Range("D36").Value = "54.46"
$ws.Range("E36").Value = "  -2.18%  "
$ws.Range("D37").Value = "0.0891"
$ws.Range("E37").Value = "  +4.42%  "
$ws.Range("D38").Value = "466.33"
$ws.Range("E38").Value = "  -1.61%  "
$ws.Range("D39").Value = "0.0412"
$ws.Range("E39").Value = "  -2.60%  "
$ws.Range("E40").Value = "  -1.63%  "
$ws.Range("D41").Value = "8.59"
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("D42").Value = "3.010.04"
$ws.Range("E42").Value = "  -3.73%  "
$ws.Range("E43").Value = "  -4.46%  "
$ws.Range("E44").Value = "  -2.90%  "
$ws.Range("E45").Value = "  -1.49%  "
$ws.Range("D46").Value = "28.01"
$ws.Range("E46").Value = "  -3.69%  "
$ws.Range("D47").Value = "0.0₃0569"
$ws.Range("E47").Value = "  +1.27%  "
$ws.Range("E49").Value = "  -2.14%  "
$ws.Range("E50").Value = "  -3.01%  "
$ws.Range("D51").Value = "117.38"
$ws.Range("E51").Value = "  -0.49%  "
